# Author's commit: "Add files via upload" — the quiz data was revised:
#  - the old row 4 ("3. 다음 문장의 뜻을 올바르게 해석한 것은? \"뭐 뭇나?\"" duplicate
#    question) was removed, shifting every later row up by one;
#  - several of the now-shifted question numbers/prompts and a few choice
#    cells were touched up to match their new position / wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Remove the old row 4 in its entirety; rows 5-13 shift up to 4-12,
#    carrying their styles/row-heights with them (matches the target
#    dimension A1:F12 and the per-row ht="33"/s="1" pattern).
$ws.Rows.Item(4).Delete()

# 2) Touch up the cells whose text differs from a pure shift-by-one.

# Row 4 (was row 5): renumber 4 -> 3, drop the leading "어어어 " and make
# "파이다." underlined via a rich-text run.
$ws.Range("A4").Value = "3. 다음 중 밑줄 친 문자의 뜻을 올바르게 해석한 것은?`n그 옷 파이다."
$ws.Range("A4").Characters(38, 4).Font.Underline = $true

# Row 5 (was row 6): renumber 5 -> 4
$ws.Range("A5").Value = "4. 다음의 뜻 풀이 중 옳지 못한 것은?"

# Row 6 (was row 7): renumber 6 -> 5; swap 행님/솔찬히 out for 솔찬히/히야
$ws.Range("A6").Value = "5. 다음 중 경상도 사투리가 아닌 것은?"
$ws.Range("D6").Value = "솔찬히"
$ws.Range("F6").Value = "히야"

# Row 7 (was row 8): renumber 7 -> 6; swap 요령껏 해라 / 제대로 해라
$ws.Range("A7").Value = "6. 다음 문자의 뜻을 올바르게 해석한 것은?`n`"단디 해라.`""
$ws.Range("C7").Value = "제대로 해라"
$ws.Range("E7").Value = "요령껏 해라"

# Row 8 (was row 9): renumber 8 -> 7
$ws.Range("A8").Value = "7. 다음 문장의 뜻을 올바르게 해석한 것은?`n이 귤 이래 세가라바가 물수있겠나"

# Row 9 (was row 10): renumber 9 -> 8
$ws.Range("A9").Value = "8. 다음 문장의 뜻을 올바르게 해석한 것은?`n아따 디다"

# Row 10 (was row 11): renumber 10 -> 9
$ws.Range("A10").Value = "9. 다음 문장의 뜻을 올바르게 해석한 것은?`n야는 와이래 분답노"

# Row 11 (was row 12): renumber 11 -> 10
$ws.Range("A11").Value = "10. 다음 문장의 뜻을 올바르게 해석한 것은?`n하늘에 별이 천지삐까리네"

# Row 12 (was row 13): renumber 12 -> 11 and drop the "혼자 먹니/혼자 묵노"
# tails from every cell in the row; C12's wording also changes.
$ws.Range("A12").Value = "11. 다음 문장의 뜻을 올바르게 해석한 것은?`n좋은거 있으면 농갈라무야지"
$ws.Range("B12").Value = "좋은거 있으면 나눠먹어야지"
$ws.Range("C12").Value = "좋은거 있으면 아껴먹어야지"
$ws.Range("D12").Value = "좋은거 있으면 먹어보라고 해야지"
$ws.Range("E12").Value = "좋은거 있으면 말을 해야지"
$ws.Range("F12").Value = "좋은거 있으면 나눠먹어야지"

# 3) Match the final selection / scroll position recorded in the sheet view.
$ws.Range("F13").Select()
